# Auto-generated edit script applying numeric corrections to the Sheets
# (ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR) per the commit diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 2092.8281  # H17: 2091.4 -> 2092.8281
$ws.Cells.Item(17, 10).Value = 2092.8281  # J17: 2091.4 -> 2092.8281
$ws.Cells.Item(17, 12).Value = 6278.4843  # L17: 6274.200000000001 -> 6278.4843
$ws.Cells.Item(17, 14).Value = -6614.4843  # N17: -6610.200000000001 -> -6614.4843
$ws.Cells.Item(93, 8).Value = 0  # H93: 54276 -> 0
$ws.Cells.Item(93, 9).Value = 0  # I93: 78000 -> 0
$ws.Cells.Item(93, 10).Value = 0  # J93: 30552 -> 0
$ws.Cells.Item(93, 11).Value = 0  # K93: 78000 -> 0
$ws.Cells.Item(93, 12).Value = 0  # L93: 30552 -> 0
$ws.Cells.Item(93, 13).ClearContents()  # M93: was -75504
$ws.Cells.Item(93, 14).ClearContents()  # N93: was -35544
$ws.Cells.Item(106, 8).Value = 2680.4102  # H106: 2655.85 -> 2680.4102
$ws.Cells.Item(106, 9).Value = 2544.875  # I106: 2519.2122 -> 2544.875
$ws.Cells.Item(106, 11).Value = 2544.875  # K106: 2519.2122 -> 2544.875
$ws.Cells.Item(106, 13).Value = -1913.875  # M106: -1888.2122 -> -1913.875
$ws.Cells.Item(132, 8).Value = 1391.4576  # H132: 1370.4746 -> 1391.4576
$ws.Cells.Item(132, 9).Value = 1350.8948  # I132: 1329.1754 -> 1350.8948
$ws.Cells.Item(132, 11).Value = 4052.6844  # K132: 3987.5262 -> 4052.6844
$ws.Cells.Item(132, 13).Value = -1522.6844  # M132: -1457.5262 -> -1522.6844
$ws.Cells.Item(138, 8).Value = 5726.655  # H138: 5690.3257 -> 5726.655
$ws.Cells.Item(138, 9).Value = 2328.5667  # I138: 2406.3794 -> 2328.5667
$ws.Cells.Item(138, 10).Value = 7614.4814  # J138: 7361.1055 -> 7614.4814
$ws.Cells.Item(138, 11).Value = 6985.7001  # K138: 7219.138199999999 -> 6985.7001
$ws.Cells.Item(138, 12).Value = 22843.4442  # L138: 22083.3165 -> 22843.4442
$ws.Cells.Item(138, 13).Value = -1845.7001  # M138: -2079.138199999999 -> -1845.7001
$ws.Cells.Item(138, 14).Value = -33123.4442  # N138: -32363.3165 -> -33123.4442
$ws.Cells.Item(141, 8).Value = 6735.476  # H141: 6433.864 -> 6735.476
$ws.Cells.Item(141, 9).Value = 6692.0557  # I141: 6345.1055 -> 6692.0557
$ws.Cells.Item(141, 11).Value = 20076.1671  # K141: 19035.3165 -> 20076.1671
$ws.Cells.Item(141, 13).Value = -14896.1671  # M141: -13855.3165 -> -14896.1671

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 3654.3372  # H32: 3624.5977 -> 3654.3372
$ws.Cells.Item(32, 9).Value = 3135.8193  # I32: 3111.1904 -> 3135.8193
$ws.Cells.Item(32, 11).Value = 3135.8193  # K32: 3111.1904 -> 3135.8193
$ws.Cells.Item(32, 13).Value = -2848.8193  # M32: -2824.1904 -> -2848.8193
$ws.Cells.Item(45, 8).Value = 2105.077  # H45: 2216.3914 -> 2105.077
$ws.Cells.Item(45, 9).Value = 1760.5454  # I45: 1840.8948 -> 1760.5454
$ws.Cells.Item(45, 11).Value = 1760.5454  # K45: 1840.8948 -> 1760.5454
$ws.Cells.Item(45, 13).Value = -1383.5454  # M45: -1463.8948 -> -1383.5454
$ws.Cells.Item(74, 8).Value = 2160.484  # H74: 2077.9697 -> 2160.484
$ws.Cells.Item(74, 9).Value = 2044.7916  # I74: 1948.9615 -> 2044.7916
$ws.Cells.Item(74, 11).Value = 2044.7916  # K74: 1948.9615 -> 2044.7916
$ws.Cells.Item(74, 13).Value = -1170.7916  # M74: -1074.9615 -> -1170.7916
$ws.Cells.Item(77, 8).Value = 2160.484  # H77: 2077.9697 -> 2160.484
$ws.Cells.Item(77, 9).Value = 2044.7916  # I77: 1948.9615 -> 2044.7916
$ws.Cells.Item(77, 11).Value = 10223.958  # K77: 9744.807499999999 -> 10223.958
$ws.Cells.Item(77, 13).Value = -5855.958000000001  # M77: -5376.807499999999 -> -5855.958000000001
$ws.Cells.Item(88, 8).Value = 7633.1665  # H88: 8798.4 -> 7633.1665
$ws.Cells.Item(88, 10).Value = 7268.6665  # J88: 9999.5 -> 7268.6665
$ws.Cells.Item(88, 12).Value = 7268.6665  # L88: 9999.5 -> 7268.6665
$ws.Cells.Item(88, 14).Value = -8080.6665  # N88: -10811.5 -> -8080.6665
$ws.Cells.Item(91, 8).Value = 7633.1665  # H91: 8798.4 -> 7633.1665
$ws.Cells.Item(91, 10).Value = 7268.6665  # J91: 9999.5 -> 7268.6665
$ws.Cells.Item(91, 12).Value = 7268.6665  # L91: 9999.5 -> 7268.6665
$ws.Cells.Item(91, 14).Value = -10076.6665  # N91: -12807.5 -> -10076.6665
$ws.Cells.Item(96, 8).Value = 61950.668  # H96: 67666.664 -> 61950.668
$ws.Cells.Item(96, 10).Value = 61950.668  # J96: 67666.664 -> 61950.668
$ws.Cells.Item(96, 12).Value = 61950.668  # L96: 67666.664 -> 61950.668
$ws.Cells.Item(96, 14).Value = -67442.66800000001  # N96: -73158.664 -> -67442.66800000001
$ws.Cells.Item(110, 8).Value = 173986.52  # H110: 157738.88 -> 173986.52
$ws.Cells.Item(110, 9).Value = 228898.55  # I110: 201512.12 -> 228898.55
$ws.Cells.Item(110, 11).Value = 228898.55  # K110: 201512.12 -> 228898.55
$ws.Cells.Item(110, 13).Value = -226853.55  # M110: -199467.12 -> -226853.55
$ws.Cells.Item(122, 8).Value = 6149.7  # H122: 6687.375 -> 6149.7
$ws.Cells.Item(122, 9).Value = 3999  # I122: 0 -> 3999
$ws.Cells.Item(122, 11).Value = 11997  # K122: 0 -> 11997
$ws.Cells.Item(122, 13).Value = -9547  # M122: None -> -9547
$ws.Cells.Item(132, 8).Value = 3318.7446  # H132: 3721.9756 -> 3318.7446
$ws.Cells.Item(132, 9).Value = 2772.318  # I132: 3121.1052 -> 2772.318
$ws.Cells.Item(132, 11).Value = 8316.954000000002  # K132: 9363.3156 -> 8316.954000000002
$ws.Cells.Item(132, 13).Value = -5786.954000000002  # M132: -6833.3156 -> -5786.954000000002

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 71433820  # H20: 125005464 -> 71433820
$ws.Cells.Item(20, 9).Value = 125004824  # I20: 250004670 -> 125004824
$ws.Cells.Item(20, 10).Value = 5820.1665  # J20: 6254 -> 5820.1665
$ws.Cells.Item(20, 11).Value = 125004824  # K20: 250004670 -> 125004824
$ws.Cells.Item(20, 12).Value = 5820.1665  # L20: 6254 -> 5820.1665
$ws.Cells.Item(20, 13).Value = -125004577  # M20: -250004423 -> -125004577
$ws.Cells.Item(20, 14).Value = -6314.1665  # N20: -6748 -> -6314.1665
$ws.Cells.Item(86, 8).Value = 947278.75  # H86: 812008.9 -> 947278.75
$ws.Cells.Item(86, 9).Value = 1419681  # I86: 1064856.8 -> 1419681
$ws.Cells.Item(86, 10).Value = 2474.3333  # J86: 2895.8 -> 2474.3333
$ws.Cells.Item(86, 11).Value = 1419681  # K86: 1064856.8 -> 1419681
$ws.Cells.Item(86, 12).Value = 2474.3333  # L86: 2895.8 -> 2474.3333
$ws.Cells.Item(86, 13).Value = -1418558  # M86: -1063733.8 -> -1418558
$ws.Cells.Item(86, 14).Value = -4720.3333  # N86: -5141.8 -> -4720.3333
$ws.Cells.Item(89, 8).Value = 947278.75  # H89: 812008.9 -> 947278.75
$ws.Cells.Item(89, 9).Value = 1419681  # I89: 1064856.8 -> 1419681
$ws.Cells.Item(89, 10).Value = 2474.3333  # J89: 2895.8 -> 2474.3333
$ws.Cells.Item(89, 11).Value = 7098405  # K89: 5324284 -> 7098405
$ws.Cells.Item(89, 12).Value = 12371.6665  # L89: 14479 -> 12371.6665
$ws.Cells.Item(89, 13).Value = -7092789  # M89: -5318668 -> -7092789
$ws.Cells.Item(89, 14).Value = -23603.6665  # N89: -25711 -> -23603.6665
$ws.Cells.Item(94, 8).Value = 1217.5714  # H94: 1107.0454 -> 1217.5714
$ws.Cells.Item(94, 9).Value = 1246.4736  # I94: 1123.45 -> 1246.4736
$ws.Cells.Item(94, 11).Value = 1246.4736  # K94: 1123.45 -> 1246.4736
$ws.Cells.Item(94, 13).Value = -795.4736  # M94: -672.45 -> -795.4736
$ws.Cells.Item(105, 8).Value = 2256.4614  # H105: 2079.484 -> 2256.4614
$ws.Cells.Item(105, 9).Value = 2206.3914  # I105: 2019.3928 -> 2206.3914
$ws.Cells.Item(105, 11).Value = 2206.3914  # K105: 2019.3928 -> 2206.3914
$ws.Cells.Item(105, 13).Value = -459.3914  # M105: -272.3928000000001 -> -459.3914
$ws.Cells.Item(107, 8).Value = 314546.3  # H107: 335497.56 -> 314546.3
$ws.Cells.Item(107, 9).Value = 1909.6072  # I107: 2035.1538 -> 1909.6072
$ws.Cells.Item(107, 11).Value = 1909.6072  # K107: 2035.1538 -> 1909.6072
$ws.Cells.Item(107, 13).Value = 10.39280000000008  # M107: -115.1538 -> 10.39280000000008
$ws.Cells.Item(134, 8).Value = 26567.232  # H134: 25976.955 -> 26567.232
$ws.Cells.Item(134, 9).Value = 3430.6667  # I134: 3359.775 -> 3430.6667
$ws.Cells.Item(134, 11).Value = 10292.0001  # K134: 10079.325 -> 10292.0001
$ws.Cells.Item(134, 13).Value = -7757.000100000001  # M134: -7544.325000000001 -> -7757.000100000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 5199.8  # H16: 3241.6667 -> 5199.8
$ws.Cells.Item(16, 9).Value = 9999  # I16: 2877.7778 -> 9999
$ws.Cells.Item(16, 10).Value = 4000  # J16: 4333.3335 -> 4000
$ws.Cells.Item(16, 11).Value = 9999  # K16: 2877.7778 -> 9999
$ws.Cells.Item(16, 12).Value = 4000  # L16: 4333.3335 -> 4000
$ws.Cells.Item(16, 13).Value = -9712  # M16: -2590.7778 -> -9712
$ws.Cells.Item(16, 14).Value = -4574  # N16: -4907.3335 -> -4574
$ws.Cells.Item(113, 8).Value = 5199.8  # H113: 3241.6667 -> 5199.8
$ws.Cells.Item(113, 9).Value = 9999  # I113: 2877.7778 -> 9999
$ws.Cells.Item(113, 10).Value = 4000  # J113: 4333.3335 -> 4000
$ws.Cells.Item(113, 11).Value = 9999  # K113: 2877.7778 -> 9999
$ws.Cells.Item(113, 12).Value = 4000  # L113: 4333.3335 -> 4000
$ws.Cells.Item(113, 13).Value = -7829  # M113: -707.7777999999998 -> -7829
$ws.Cells.Item(113, 14).Value = -8340  # N113: -8673.333500000001 -> -8340
$ws.Cells.Item(132, 8).Value = 1883.6487  # H132: 1852.5 -> 1883.6487
$ws.Cells.Item(132, 9).Value = 1461.6177  # I132: 1439.8572 -> 1461.6177
$ws.Cells.Item(132, 11).Value = 4384.8531  # K132: 4319.571599999999 -> 4384.8531
$ws.Cells.Item(132, 13).Value = -1854.8531  # M132: -1789.571599999999 -> -1854.8531
$ws.Cells.Item(134, 8).Value = 359893.34  # H134: 359895.16 -> 359893.34
$ws.Cells.Item(134, 9).Value = 2985.1667  # I134: 2987.25 -> 2985.1667
$ws.Cells.Item(134, 11).Value = 8955.500100000001  # K134: 8961.75 -> 8955.500100000001
$ws.Cells.Item(134, 13).Value = -6420.500100000001  # M134: -6426.75 -> -6420.500100000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(132, 8).Value = 1005539.06  # H132: 1105544 -> 1005539.06
$ws.Cells.Item(132, 9).Value = 205186.6  # I132: 255110.75 -> 205186.6
$ws.Cells.Item(132, 11).Value = 1846679.4  # K132: 2295996.75 -> 1846679.4
$ws.Cells.Item(132, 13).Value = -1844149.4  # M132: -2293466.75 -> -1844149.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 205.7037  # H2: 217.57692 -> 205.7037
$ws.Cells.Item(2, 9).Value = 212  # I2: 217.13637 -> 212
$ws.Cells.Item(2, 10).Value = 178  # J2: 220 -> 178
$ws.Cells.Item(2, 11).Value = 212  # K2: 217.13637 -> 212
$ws.Cells.Item(2, 12).Value = 178  # L2: 220 -> 178
$ws.Cells.Item(2, 13).Value = -99  # M2: -104.13637 -> -99
$ws.Cells.Item(2, 14).Value = -404  # N2: -446 -> -404
$ws.Cells.Item(80, 8).Value = 913494.6  # H80: 744828.4 -> 913494.6
$ws.Cells.Item(80, 9).Value = 774315.25  # I80: 592770.0600000001 -> 774315.25
$ws.Cells.Item(80, 10).Value = 1114531.5  # J80: 1003327.5 -> 1114531.5
$ws.Cells.Item(80, 11).Value = 774315.25  # K80: 592770.0600000001 -> 774315.25
$ws.Cells.Item(80, 12).Value = 1114531.5  # L80: 1003327.5 -> 1114531.5
$ws.Cells.Item(80, 13).Value = -773317.25  # M80: -591772.0600000001 -> -773317.25
$ws.Cells.Item(80, 14).Value = -1116527.5  # N80: -1005323.5 -> -1116527.5
$ws.Cells.Item(83, 8).Value = 913494.6  # H83: 744828.4 -> 913494.6
$ws.Cells.Item(83, 9).Value = 774315.25  # I83: 592770.0600000001 -> 774315.25
$ws.Cells.Item(83, 10).Value = 1114531.5  # J83: 1003327.5 -> 1114531.5
$ws.Cells.Item(83, 11).Value = 3871576.25  # K83: 2963850.3 -> 3871576.25
$ws.Cells.Item(83, 12).Value = 5572657.5  # L83: 5016637.5 -> 5572657.5
$ws.Cells.Item(83, 13).Value = -3866584.25  # M83: -2958858.3 -> -3866584.25
$ws.Cells.Item(83, 14).Value = -5582641.5  # N83: -5026621.5 -> -5582641.5
$ws.Cells.Item(97, 8).Value = 594.0833  # H97: 610.2174 -> 594.0833
$ws.Cells.Item(97, 9).Value = 648.7  # I97: 671.1053000000001 -> 648.7
$ws.Cells.Item(97, 11).Value = 648.7  # K97: 671.1053000000001 -> 648.7
$ws.Cells.Item(97, 13).Value = -152.7  # M97: -175.1053000000001 -> -152.7
$ws.Cells.Item(122, 8).Value = 2921.75  # H122: 2999.2812 -> 2921.75
$ws.Cells.Item(122, 9).Value = 2580.7083  # I122: 2645.1738 -> 2580.7083
$ws.Cells.Item(122, 10).Value = 3944.875  # J122: 3904.2222 -> 3944.875
$ws.Cells.Item(122, 11).Value = 7742.124899999999  # K122: 7935.5214 -> 7742.124899999999
$ws.Cells.Item(122, 12).Value = 11834.625  # L122: 11712.6666 -> 11834.625
$ws.Cells.Item(122, 13).Value = -5292.124899999999  # M122: -5485.5214 -> -5292.124899999999
$ws.Cells.Item(122, 14).Value = -16734.625  # N122: -16612.6666 -> -16734.625
$ws.Cells.Item(132, 8).Value = 22733.18  # H132: 22734.26 -> 22733.18
$ws.Cells.Item(132, 9).Value = 2944.6445  # I132: 2945.8445 -> 2944.6445
$ws.Cells.Item(132, 11).Value = 8833.933499999999  # K132: 8837.533500000001 -> 8833.933499999999
$ws.Cells.Item(132, 13).Value = -6303.933499999999  # M132: -6307.533500000001 -> -6303.933499999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 998  # H22: 600 -> 998
$ws.Cells.Item(22, 9).Value = 0  # I22: 600 -> 0
$ws.Cells.Item(22, 10).Value = 998  # J22: 0 -> 998
$ws.Cells.Item(22, 11).Value = 0  # K22: 600 -> 0
$ws.Cells.Item(22, 12).Value = 998  # L22: 0 -> 998
$ws.Cells.Item(22, 13).ClearContents()  # M22: was -305
$ws.Cells.Item(22, 14).Value = -1588  # N22: None -> -1588
$ws.Cells.Item(27, 8).Value = 998  # H27: 600 -> 998
$ws.Cells.Item(27, 9).Value = 0  # I27: 600 -> 0
$ws.Cells.Item(27, 10).Value = 998  # J27: 0 -> 998
$ws.Cells.Item(27, 11).Value = 0  # K27: 600 -> 0
$ws.Cells.Item(27, 12).Value = 998  # L27: 0 -> 998
$ws.Cells.Item(27, 13).ClearContents()  # M27: was -493
$ws.Cells.Item(27, 14).Value = -1212  # N27: None -> -1212
$ws.Cells.Item(93, 8).Value = 3343.9211  # H93: 3343.9473 -> 3343.9211
$ws.Cells.Item(93, 9).Value = 3232.64  # I93: 3232.68 -> 3232.64
$ws.Cells.Item(93, 11).Value = 3232.64  # K93: 3232.68 -> 3232.64
$ws.Cells.Item(93, 13).Value = -1984.64  # M93: -1984.68 -> -1984.64
$ws.Cells.Item(112, 8).Value = 100000  # H112: 0 -> 100000
$ws.Cells.Item(112, 10).Value = 100000  # J112: 0 -> 100000
$ws.Cells.Item(112, 12).Value = 100000  # L112: 0 -> 100000
$ws.Cells.Item(112, 14).Value = -102954  # N112: None -> -102954
$ws.Cells.Item(133, 8).Value = 58299.6  # H133: 64166.5 -> 58299.6
$ws.Cells.Item(133, 10).Value = 58299.6  # J133: 64166.5 -> 58299.6
$ws.Cells.Item(133, 12).Value = 58299.6  # L133: 64166.5 -> 58299.6
$ws.Cells.Item(133, 14).Value = -63359.6  # N133: -69226.5 -> -63359.6
$ws.Cells.Item(136, 8).Value = 381447.94  # H136: 355693 -> 381447.94
$ws.Cells.Item(136, 9).Value = 676273.0600000001  # I136: 597652.9 -> 676273.0600000001
$ws.Cells.Item(136, 11).Value = 2028819.18  # K136: 1792958.7 -> 2028819.18
$ws.Cells.Item(136, 13).Value = -2026269.18  # M136: -1790408.7 -> -2026269.18

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 15125  # H81: 10195.467 -> 15125
$ws.Cells.Item(81, 9).Value = 3874.5  # I81: 2548.111 -> 3874.5
$ws.Cells.Item(81, 10).Value = 20750.25  # J81: 21666.5 -> 20750.25
$ws.Cells.Item(81, 11).Value = 7749  # K81: 5096.222 -> 7749
$ws.Cells.Item(81, 12).Value = 41500.5  # L81: 43333 -> 41500.5
$ws.Cells.Item(81, 13).Value = -6688  # M81: -4035.222 -> -6688
$ws.Cells.Item(81, 14).Value = -43622.5  # N81: -45455 -> -43622.5
$ws.Cells.Item(84, 8).Value = 15125  # H84: 10195.467 -> 15125
$ws.Cells.Item(84, 9).Value = 3874.5  # I84: 2548.111 -> 3874.5
$ws.Cells.Item(84, 10).Value = 20750.25  # J84: 21666.5 -> 20750.25
$ws.Cells.Item(84, 11).Value = 38745  # K84: 25481.11 -> 38745
$ws.Cells.Item(84, 12).Value = 207502.5  # L84: 216665 -> 207502.5
$ws.Cells.Item(84, 13).Value = -33441  # M84: -20177.11 -> -33441
$ws.Cells.Item(84, 14).Value = -218110.5  # N84: -227273 -> -218110.5
$ws.Cells.Item(97, 8).Value = 49972  # H97: 0 -> 49972
$ws.Cells.Item(97, 10).Value = 49972  # J97: 0 -> 49972
$ws.Cells.Item(97, 12).Value = 49972  # L97: 0 -> 49972
$ws.Cells.Item(97, 14).Value = -51954  # N97: None -> -51954
$ws.Cells.Item(126, 8).Value = 1593  # H126: 1658.9584 -> 1593
$ws.Cells.Item(126, 9).Value = 1588.8334  # I126: 1663.619 -> 1588.8334
$ws.Cells.Item(126, 11).Value = 4766.5002  # K126: 4990.857 -> 4766.5002
$ws.Cells.Item(126, 13).Value = -2296.5002  # M126: -2520.857 -> -2296.5002
$ws.Cells.Item(132, 8).Value = 42446.44  # H132: 40848.46 -> 42446.44
$ws.Cells.Item(132, 9).Value = 1542.3889  # I132: 1508.5264 -> 1542.3889
$ws.Cells.Item(132, 11).Value = 4627.1667  # K132: 4525.5792 -> 4627.1667
$ws.Cells.Item(132, 13).Value = -2097.1667  # M132: -1995.5792 -> -2097.1667
$ws.Cells.Item(136, 8).Value = 282649.25  # H136: 282775.97 -> 282649.25
$ws.Cells.Item(136, 9).Value = 247217.53  # I136: 247350.42 -> 247217.53
$ws.Cells.Item(136, 11).Value = 741652.59  # K136: 742051.26 -> 741652.59
$ws.Cells.Item(136, 13).Value = -739102.59  # M136: -739501.26 -> -739102.59
